$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 39
    3  = 23
    4  = 58
    5  = 27
    6  = 31
    7  = 19
    8  = 43
    9  = 59
    10 = 87
    11 = 18
    12 = 60
    13 = 73
    14 = 85
    15 = 15
    16 = 97
    17 = 61
    18 = 75
    19 = 32
    20 = 33
    21 = 69
    23 = 28
    24 = 41
    25 = 45
    27 = 79
    28 = 40
    29 = 1
    30 = 42
    31 = 46
    32 = 91
    33 = 56
    34 = 77
    35 = 47
    36 = 121
    37 = 104
    38 = 29
    39 = 181
    40 = 24
    41 = 88
    42 = 199
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row]
}

$wb.Save()
